# Apply the edit described by the commit: mark requirements R18 and R19
# ("Guardar/Cargar la información actual de la partida") as done ("x")
# on the "Requerimientos" sheet, and leave that sheet as the active tab
# with the selection on B20 (instead of "Must Have" being active on D8).

$wb = $excel.ActiveWorkbook

$wsReq = $wb.Worksheets.Item("Requerimientos")

# Update the status flags for R18 and R19 from "p" (pending) to "x" (done)
$wsReq.Range("B18").Value = "x"
$wsReq.Range("B19").Value = "x"

# Make "Requerimientos" the active sheet/tab and set its selection to B20
$wsReq.Activate()
$wsReq.Range("B20").Select()

# "Must Have" should no longer be the selected tab; selecting another
# sheet above already clears its tabSelected state.
